# The codeforiati:category-name / codeforiati:group-name columns (D/E) were
# swapped, as were the codeforiati:group-code / codeforiati:category-code
# columns (F/G). This affects the header row and every data row.
#
# We use Range.Copy (rather than reading/writing .Value arrays) so that the
# text representation of numeric-looking codes (e.g. "110") is preserved
# exactly as shared-string text instead of being coerced into numbers, and
# so that no cell styles are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Scratch columns placed well beyond the used range so they never collide
# with real data, and are cleared again before we are done.
$scratch1 = "AZ"
$scratch2 = "BA"

# Swap column D (codeforiati:category-name) <-> column E (codeforiati:group-name)
$ws.Range("D1:D$lastRow").Copy($ws.Range($scratch1 + "1"))
$ws.Range("E1:E$lastRow").Copy($ws.Range("D1"))
$ws.Range($scratch1 + "1:" + $scratch1 + "$lastRow").Copy($ws.Range("E1"))

# Swap column F (codeforiati:group-code) <-> column G (codeforiati:category-code)
$ws.Range("F1:F$lastRow").Copy($ws.Range($scratch2 + "1"))
$ws.Range("G1:G$lastRow").Copy($ws.Range("F1"))
$ws.Range($scratch2 + "1:" + $scratch2 + "$lastRow").Copy($ws.Range("G1"))

# Clean up the scratch columns so the saved workbook only contains A:G.
$ws.Range($scratch1 + "1:" + $scratch2 + "$lastRow").Clear()

$excel.CutCopyMode = 0
